$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - "¿Qué se hizo ayer?" for Paula Andrea Taborda Jaramillo
$ws.Range("C10").Value = "Nada "
$ws.Range("D10").Value = "Nada"
$ws.Range("E10").Value = "Nada"
$ws.Range("F10").Value = "Ordenar la reunión de hoy"
$ws.Range("G10").Value = "Se logró el objetivo y se encontró el diagrama que faltaba"

# Row 11 - "¿Qué se hará hoy?"
$ws.Range("C11").Value = "Nada "
$ws.Range("D11").Value = "Nada"
$ws.Range("E11").Value = "Planeamos reunión para mañana"
$ws.Range("F11").Value = "Reunirnos, revisar diagramas y Mariana comenzará con las tablas"
$ws.Range("G11").Value = "Reunión para revisar y opinar el trabajo de Mariana"

# Row 12 - "¿Qué cosas se oponen?"
$ws.Range("C12").Value = "Problemas personales"
$ws.Range("D12").Value = "No hubo clase "
$ws.Range("E12").Value = "Otras materias"
$ws.Range("F12").Value = "Ninguna"
$ws.Range("G12").Value = "Parcial mañana"

# Update selection to match the saved workbook view state
$ws.Range("E13").Select() | Out-Null
